$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (duplicate "Ad Board" entry), which pushes the
# existing rows 13-16 (Tip Keg, Pager, Radio, House) down to 14-17.
$ws.Rows("13:13").Insert()

$ws.Range("A13").Value = "Ad Board"
$ws.Range("E13").Value = 1

# Append a brand new trailing row (18) for the new "What are you drinking" item.
$ws.Range("A18").Value = "What are you drinking"
$ws.Range("E18").Value = 1

# Header for the new "Export" column.
$ws.Range("F1").Value = "Export"

# Column F: Export y/n flag, Column G: free-text note (when present).
$ws.Range("F2").Value = "N"
$ws.Range("G2").Value = "(mega prims)"

$ws.Range("F3").Value = "y"
$ws.Range("G3").Value = "(some textures fail)"

$ws.Range("F4").Value = "n"
$ws.Range("G4").Value = "(too hard to get too)"

$ws.Range("F5").Value = "Y"

$ws.Range("F6").Value = "N"
$ws.Range("G6").Value = "permissions!"

$ws.Range("F7").Value = "Y"

$ws.Range("F8").Value = "y"
$ws.Range("G8").Value = "textures"

$ws.Range("F9").Value = "N"

$ws.Range("F10").Value = "Y"

$ws.Range("F11").Value = "y"
$ws.Range("G11").Value = "textures"

$ws.Range("F12").Value = "Y"

$ws.Range("F13").Formula = "=NA()"
$ws.Range("G13").Value = "copy of other"

$ws.Range("F14").Value = "Y"

$ws.Range("F15").Value = "y"
$ws.Range("G15").Value = "textures"

$ws.Range("F16").Value = "N"

$ws.Range("F18").Value = "y"
$ws.Range("G18").Value = "textures"

# Update dimension/selection ranges to match new sheet extent.
$ws.Range("E2:E18").Select()
